# Insert a new event row before row 31, shifting the existing rows 31-57
# down to 32-58, then populate the new row 31 with the "Da suddite a
# cittadine..." event data (per the commit's "latest data" refresh).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 31 - this shifts rows 31..57 down to 32..58
# and extends the used range to A1:AA58.
$ws.Rows("31").Insert()

# Populate the new row 31.
$ws.Cells.Item(31, 1).Value  = "Conferenze, Seminari, Incontri e Lezioni"
$ws.Cells.Item(31, 2).Value  = "Modena"
$ws.Cells.Item(31, 3).Value  = "Strada Vaciglio Nord, 6"
$ws.Cells.Item(31, 4).Value  = "2022-06-04T08:30:34+00:00"
$ws.Cells.Item(31, 6).Value  = "2022-06-04T08:30:59+00:00"
$ws.Cells.Item(31, 8).Value  = "2022-06-10T08:00:00+00:00"
$ws.Cells.Item(31, 9).Value  = "2022-06-10T09:00:00+00:00"
$ws.Cells.Item(31, 10).Value = "https://www.comune.modena.it/api/novita/eventi/2022/da-suddite-a-cittadine-gabriella-degli-esposti-e-le-partigiane-modenesi-nelle-fonti-documentarie-e-nelle-memorie/@@images/9a9f2a19-02f6-4cd6-a272-5cac5ae7adbd.jpeg"
$ws.Cells.Item(31, 12).Value = "2022-06-04T08:52:51+00:00"
$ws.Cells.Item(31, 13).Value = "Sala Renata Bergonzoni della Casa delle Donne"
$ws.Cells.Item(31, 14).Value = " ore 18.30"
$ws.Cells.Item(31, 16).Value = " ingresso libero"
$ws.Cells.Item(31, 19).Value = "Da suddite a cittadine. Gabriella Degli Esposti e le partigiane modenesi nelle fonti documentarie e nelle memorie"
$ws.Cells.Item(31, 22).Value = $false
$ws.Cells.Item(31, 23).Value = 41123
$ws.Cells.Item(31, 24).Value = "https://www.comune.modena.it/novita/eventi/2022/da-suddite-a-cittadine-gabriella-degli-esposti-e-le-partigiane-modenesi-nelle-fonti-documentarie-e-nelle-memorie"
$ws.Cells.Item(31, 25).Value = "44,64582"
$ws.Cells.Item(31, 26).Value = "10,92572"
$ws.Cells.Item(31, 27).Value = "POINT (10.92572 44.64582)"
